$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 39.06069033333333
$ws.Range("H2").Value = 117.182071
$ws.Range("I2").Value = 0.3672373747374215
$ws.Range("J2").Value = 0.3672373747374215
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 1.534251835372889
$ws.Range("R2").Value = 13.808266518356
$ws.Range("S2").Value = 0.2729501093695562
$ws.Range("T2").Value = 0.2729501093695562
$ws.Range("G3").Value = 39.06069033333333
$ws.Range("H3").Value = 117.182071
$ws.Range("I3").Value = 0.3672373747374215
$ws.Range("J3").Value = 0.3672373747374215
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 0.5299884666727778
$ws.Range("R3").Value = 4.769896200055
$ws.Range("S3").Value = 0.09428726536786536
$ws.Range("T3").Value = 0.09428726536786536
$ws.Range("I4").Value = 0.1367778286588004
$ws.Range("J4").Value = 0.1367778286588004
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 0.5714332175697777
$ws.Range("R4").Value = 5.142898958128
$ws.Range("S4").Value = 0.1016604677518017
$ws.Range("T4").Value = 0.1016604677518018
$ws.Range("I5").Value = 0.1367778286588004
$ws.Range("J5").Value = 0.1367778286588004
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 0.1973945918155556
$ws.Range("R5").Value = 1.77655132634
$ws.Range("S5").Value = 0.03511736090699862
$ws.Range("T5").Value = 0.03511736090699862
$ws.Range("G6").Value = 28.68702533333333
$ws.Range("H6").Value = 86.061076
$ws.Range("I6").Value = 0.2697071603839269
$ws.Range("J6").Value = 0.2697071603839269
$ws.Range("M6").Value = 0.03927866666666666
$ws.Range("N6").Value = 0.117836
$ws.Range("O6").Value = 0.7432525340448212
$ws.Range("P6").Value = 0.7432525340448213
$ws.Range("Q6").Value = 1.126788105726222
$ws.Range("R6").Value = 10.141092951536
$ws.Range("S6").Value = 0.2004605304053867
$ws.Range("T6").Value = 0.2004605304053867
$ws.Range("G7").Value = 28.68702533333333
$ws.Range("H7").Value = 86.061076
$ws.Range("I7").Value = 0.2697071603839269
$ws.Range("J7").Value = 0.2697071603839269
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01356833333333333
$ws.Range("N7").Value = 0.040705
$ws.Range("O7").Value = 0.2567474659551788
$ws.Range("P7").Value = 0.2567474659551788
$ws.Range("Q7").Value = 0.3892351220644444
$ws.Range("R7").Value = 3.50311609858
$ws.Range("S7").Value = 0.06924662997854023
$ws.Range("T7").Value = 0.06924662997854021
$ws.Range("G8").Value = 24.06770466666667
$ws.Range("H8").Value = 72.203114
$ws.Range("I8").Value = 0.2262776362198511
$ws.Range("J8").Value = 0.2262776362198511
$ws.Range("M8").Value = 0.03927866666666666
$ws.Range("N8").Value = 0.117836
$ws.Range("O8").Value = 0.7432525340448212
$ws.Range("P8").Value = 0.7432525340448213
$ws.Range("Q8").Value = 0.9453473490337777
$ws.Range("R8").Value = 8.508126141304
$ws.Range("S8").Value = 0.1681814265180765
$ws.Range("T8").Value = 0.1681814265180766
$ws.Range("G9").Value = 24.06770466666667
$ws.Range("H9").Value = 72.203114
$ws.Range("I9").Value = 0.2262776362198511
$ws.Range("J9").Value = 0.2262776362198511
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01356833333333333
$ws.Range("N9").Value = 0.040705
$ws.Range("O9").Value = 0.2567474659551788
$ws.Range("P9").Value = 0.2567474659551788
$ws.Range("Q9").Value = 0.3265586394855556
$ws.Range("R9").Value = 2.93902775537
$ws.Range("S9").Value = 0.05809620970177455
$ws.Range("T9").Value = 0.05809620970177455
